$d = $word.ActiveDocument

$pairs = @(
    @("561÷2=280, 1", "501÷4=125, 1"),
    @("828÷2=414, 0", "616÷8=77, 0"),
    @("292÷5=58, 2", "687÷8=85, 7"),
    @("751÷7=107, 2", "622÷8=77, 6"),
    @("799÷2=399, 1", "543÷2=271, 1"),
    @("572÷4=143, 0", "306÷4=76, 2"),
    @("733÷9=81, 4", "629÷3=209, 2"),
    @("540÷2=270, 0", "974÷6=162, 2"),
    @("612÷9=68, 0", "316÷5=63, 1"),
    @("566÷6=94, 2", "653÷4=163, 1"),
    @("151÷2=75, 1", "610÷6=101, 4"),
    @("356÷3=118, 2", "362÷9=40, 2"),
    @("951÷8=118, 7", "796÷8=99, 4"),
    @("241÷5=48, 1", "272÷6=45, 2"),
    @("212÷7=30, 2", "205÷5=41, 0"),
    @("501÷6=83, 3", "615÷5=123, 0"),
    @("474÷3=158, 0", "353÷4=88, 1"),
    @("586÷3=195, 1", "757÷6=126, 1"),
    @("600÷9=66, 6", "971÷7=138, 5"),
    @("305÷6=50, 5", "410÷2=205, 0"),
    @("346÷3=115, 1", "866÷5=173, 1"),
    @("933÷4=233, 1", "794÷7=113, 3"),
    @("180÷2=90, 0", "163÷7=23, 2"),
    @("180÷6=30, 0", "831÷4=207, 3"),
    @("662÷9=73, 5", "541÷8=67, 5")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
